$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "26.431.48", "4.00") that must stay
# text, not be auto-coerced to numbers by Excel - so force text format first.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.428.70"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "1.604.29"
$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "212.29"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -0.33%  "

$ws.Range("D9").Value = "0.0607"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").Value = "19.34"
$ws.Range("E10").Value = "  +1.49%  "

$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").Value = "1.829.56"
$ws.Range("E12").Value = "  +0.82%  "

$ws.Range("D13").Value = "1.611.37"
$ws.Range("E13").Value = "  +1.25%  "

$ws.Range("D14").Value = "4.01"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").Value = "0.508"
$ws.Range("E15").Value = "  -0.47%  "

$ws.Range("D16").Value = "63.72"
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").Value = "233.91"
$ws.Range("E17").Value = "  +8.20%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "26.418.96"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").Value = "7.73"
$ws.Range("E19").Value = "  +5.82%  "

$ws.Range("D20").Value = "0.0₃0725"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "8.98"
$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  +0.59%  "

$ws.Range("D25").Value = "147.17"
$ws.Range("E25").Value = "  +1.84%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").Value = "7.00"
$ws.Range("E27").Value = "  +0.34%  "

$ws.Range("E28").Value = "  +1.40%  "

$ws.Range("D29").Value = "15.49"
$ws.Range("E29").Value = "  +2.36%  "

$ws.Range("D30").Value = "0.0496"
$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").Value = "1.503.81"
$ws.Range("E32").Value = "  +5.19%  "

$ws.Range("D33").Value = "3.23"
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").Value = "2.95"
$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("E35").Value = "  -0.30%  "

$ws.Range("D36").Value = "1.48"
$ws.Range("E36").Value = "  +0.87%  "

$ws.Range("D37").Value = "0.569"
$ws.Range("E37").Value = "  -2.51%  "

$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("D39").Value = "0.824"
$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("D40").Value = "5.81"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "2.18"
$ws.Range("E42").Value = "  +2.13%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "0.946"
$ws.Range("E43").Value = "  -4.52%  "

$ws.Range("D44").Value = "1.742.18"
$ws.Range("E44").Value = "  +0.86%  "

$ws.Range("D45").Value = "0.763"
$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("D46").Value = "60.99"
$ws.Range("E46").Value = "  -0.20%  "

$ws.Range("D47").Value = "89.14"
$ws.Range("E47").Value = "  +2.58%  "

$ws.Range("D48").Value = "1.50"
$ws.Range("E48").Value = "  +0.90%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").Value = "0.0965"
$ws.Range("E50").Value = "  +1.23%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.48"
$ws.Range("E51").Value = "  +1.58%  "
